$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the dataset. This pushes the
# existing rows 35-104 down to 36-105 (dimension grows from R104 to R105),
# and the new row 35 is populated with the new record's data.
$ws.Rows(35).Insert()

$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value = 44775
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = 100112001
$ws.Range("G35").Value = "Berenjena"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 150
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 11333
$ws.Range("N35").Value = "$/caja 60 unidades"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 189
$ws.Range("Q35").Value = 60
$ws.Range("R35").Value = "Hortaliza"
